# Fix mojibake "Â±" -> "±" in the results table (columns B, C, D; rows 2-17)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$badChar = [string][char]0x00C2 + [string][char]0x00B1  # "Â±"
$goodChar = [string][char]0x00B1                         # "±"

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$row")
        $val = $cell.Value2
        if ($val -ne $null -and $val.Contains($badChar)) {
            $cell.Value2 = $val.Replace($badChar, $goodChar)
        }
    }
}
